# Update countries & provincias Spain
# Applies the 26-Mar-2020 05:42 -> 06:12 data refresh:
#   - numeric updates for several country rows (new case counts)
#   - a handful of countries (Mexico, Kazajistan) overtook their neighbours
#     in the case-count ranking, so those rows swap places with the ones
#     they passed
#   - the "last updated" timestamp footer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric refreshes (country stays on its row) ---------------

# Estados Unidos (row 6)
$ws.Range("B6").Value = 68472
$ws.Range("C6").Value = 261
$ws.Range("E6").Value = 67046

# Australia (row 20)
$ws.Range("D20").Value = 170
$ws.Range("E20").Value = 2546

# Pakistan (row 33)
$ws.Range("B33").Value = 1093
$ws.Range("C33").Value = 30
$ws.Range("E33").Value = 1064

# India (row 44)
$ws.Range("B44").Value = 665
$ws.Range("C44").Value = 8
$ws.Range("E44").Value = 610

# --- Mexico overtakes Colombia / Egipto / Croacia / Barein / Hong Kong -

$ws.Range("A53").Value = "Mexico"
$ws.Range("B53").Value = 475
$ws.Range("C53").Value = 70
$ws.Range("D53").Value = 4
$ws.Range("E53").Value = 465
$ws.Range("F53").Value = 1
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 6

$ws.Range("A54").Value = "Colombia"
$ws.Range("B54").Value = 470
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 8
$ws.Range("E54").Value = 458
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 4

$ws.Range("A55").Value = "Egipto"
$ws.Range("B55").Value = 456
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 95
$ws.Range("E55").Value = 340
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 21

$ws.Range("A56").Value = "Croacia"
$ws.Range("B56").Value = 442
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 22
$ws.Range("E56").Value = 419
$ws.Range("F56").Value = 6
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 1

$ws.Range("A57").Value = "Barein"
$ws.Range("B57").Value = 419
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 177
$ws.Range("E57").Value = 238
$ws.Range("F57").Value = 2
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 4

$ws.Range("A58").Value = "Hong Kong"
$ws.Range("B58").Value = 411
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 102
$ws.Range("E58").Value = 305
$ws.Range("F58").Value = 4
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 4

# --- Kazajistan overtakes Bielorrusia / Afganistan ----------------------

$ws.Range("A100").Value = "Kazajistan"
$ws.Range("B100").Value = 88
$ws.Range("C100").Value = 7
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 88
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 0

$ws.Range("A101").Value = "Bielorrusia"
$ws.Range("B101").Value = 86
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 29
$ws.Range("E101").Value = 57
$ws.Range("F101").Value = 2
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0

$ws.Range("A102").Value = "Afganistan"
$ws.Range("B102").Value = 84
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 2
$ws.Range("E102").Value = 80
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 2

# --- Footer timestamp ---------------------------------------------------

$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 06:12"
